# "New cases 25 Nov" - add new PlaceHolder/MetaData test cases to Test_Data.xlsx
# and tweak the Parameters text of an existing PlaceHolder test case (row 55).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# 1. Row 55 (existing "YES,PlaceHolder,Test Publish,TESTUSN" case) gets a
#    trailing "~" appended onto "Test Publish" in its Parameters column.
# ---------------------------------------------------------------------------
$ws.Cells.Item(55, 4).Value = "YES,PlaceHolder,Test Publish~,TESTUSN"

# ---------------------------------------------------------------------------
# 2. Append seven new rows (105-111) to the "Table2" ListObject, carrying the
#    formatting of the most similar existing rows so the new cells pick up
#    matching styles/borders.
# ---------------------------------------------------------------------------

function Add-DataRow($Row, $A, $B, $C, $D, $E, $F, $DStyleSrc, $FStyleSrc) {
    $lo.ListRows.Add() | Out-Null

    # Copy the bulk A:E formatting from row 55 (No / TC_ID / ScriptName /
    # Parameters / MetaData - the same layout the new rows use).
    $ws.Range("A55:E55").Copy() | Out-Null
    $target = "A" + $Row + ":E" + $Row
    $ws.Range($target).PasteSpecial(-4122) | Out-Null

    # D column sometimes needs the alternate border style used at D2/D109.
    $ws.Range($DStyleSrc).Copy() | Out-Null
    $ws.Range("D" + $Row).PasteSpecial(-4122) | Out-Null

    # F (Description) column formatting.
    $ws.Range($FStyleSrc).Copy() | Out-Null
    $ws.Range("F" + $Row).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($Row, 1).Value = $A
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
}

Add-DataRow 105 "No" "FW_UI_0098" "Verify_Alert_Publish" `
    "YES,PlaceHolder,~company~ ~ric~ - SALES AND REVENUES WERE UP `$1.780 BLN FROM Q1 OF 2007. SALES VOLUME IMPROVED `$1.087,TESTUSN" `
    "MetaData" `
    "Aim of the Script is to verify that user should not be able to publish an alert with incomplete placeholders (with incomplete company and RIC), so that I do not publish out an incomplete alert " `
    "D55" "F2"

Add-DataRow 106 "No" "FW_UI_0099" "Verify_Alert_Publish" `
    "YES,PlaceHolder,Company ~ric~ - SALES AND REVENUES WERE UP `$1.780 BLN FROM Q1 OF 2007. SALES VOLUME IMPROVED,TESTUSN" `
    "MetaData" `
    "Aim of the Script is to verify that user should not be able to publish an alert with incomplete placeholders (with an incomplete RIC), so that I do not publish out an incomplete alert " `
    "D55" "F2"

Add-DataRow 107 "No" "FW_UI_0100" "Verify_Alert_Publish" `
    "YES,PlaceHolder,~company~ ric - SALES AND REVENUES WERE UP `$1.780 BLN FROM Q1 OF 2018. SALES VOLUME IMPROVED,TESTUSN" `
    "MetaData" `
    "Aim of the Script is to verify that user should not be able to publish an alert with incomplete placeholders (with an incomplete Company), so that I do not publish out an incomplete alert " `
    "D55" "F2"

Add-DataRow 108 "No" "FW_UI_0101" "Verify_MetaData_Inputs" `
    "Products,SCAN;UKP;SUDB;HX;SUKP,Publish" `
    "MetaData" `
    "Aim of the Script is to verify whether user is able to add multiple product codes in alert editor and publish" `
    "D55" "F57"

Add-DataRow 109 "No" "FW_UI_0102" "Verify_MetaData_Inputs" `
    "Topics,SASIAE;SANPRO;HAND;HARW;HACK,Publish" `
    "MetaData" `
    "Aim of the Script is to verify whether user is able to add multiple topic codes in alert editor and publish" `
    "D2" "F57"

Add-DataRow 110 "Yes" "FW_UI_0103" "Verify_MetaData_Inputs" `
    "RICS,H.N;D11.HN,Publish" `
    "MetaData" `
    "Aim of the Script is to verify whether user is able to add multiple RICs in alert editor and publish" `
    "D2" "F57"

Add-DataRow 111 "No" "FW_UI_0104" "Verify_MetaData_Inputs" `
    "NamedItems,ABS/;AB/CN,Publish" `
    "MetaData" `
    "Aim of the Script is to verify whether user is able to add multiple named items in alert editor and publish" `
    "D2" "F57"

# ---------------------------------------------------------------------------
# 3. Leave the grid focused/scrolled roughly where the author left it.
# ---------------------------------------------------------------------------
$ws.Range("D110").Select() | Out-Null

Write-Output "Added rows 105-111 and updated row 55 Parameters."
